# Updated detailed indicator quantile results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Git Commit ID recorded for the IndicatorQuantiles.R script
# (column AJ, "ScriptLatestRunVersion") for all data rows (2-80).
$ws.Range("AJ2:AJ80").Value = "IndicatorQuantiles.R, Git Commit ID: db49f0f869e1f5a8558dc746458075a467cf2c41"

# Update the recorded process id (column AH, "pid") for all data rows (2-80).
$ws.Range("AH2:AH80").Value = 21528
